$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 0
$ws.Range("E6").Value = 191.9

[void]$ws.Range("E2").Select()
